$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Problem Premise:" -> "Problem Statement:"
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Premise", $false, $false, $false, $false, $false, $true, 1, $false, "Statement", 2) | Out-Null

# ------------------------------------------------------------------
# 2. "Attached is the road map ..." -> "There is a road map ..."
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Attached is the", $false, $false, $false, $false, $false, $true, 1, $false, "There is a", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Re-colour the sentence "if light is green on A-B path, vehicles
#    can go in the direction A-B and B-A not across (E-D or D-E)."
#    from the inherited red to automatic/Text 1, leaving the leading
#    ". " in its original red colour.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("if light is green on A-B path, vehicles can go in the direction A-B and B-A not across (E-D or D-E).")
if ($found) {
    $rng.Font.TextColor.RGB = 0
    $rng.Font.TextColor.ObjectThemeColor = 13
}

# ------------------------------------------------------------------
# 4. Drop a "_GoBack" bookmark in the middle of "standing" (st|anding),
#    mirroring the last-edit marker Word leaves behind.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("how many vehicles are st")
if ($found) {
    $pos = $rng.End
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

# ------------------------------------------------------------------
# 5. Tag the image run with the en-US language (was eastAsia en-IN).
# ------------------------------------------------------------------
if ($d.InlineShapes.Count -ge 1) {
    $d.InlineShapes.Item(1).Range.LanguageID = "en-US"
}

# ------------------------------------------------------------------
# 6. Remove the trailing "How to use / Example / (blank) / Maven
#    Project / (blank)" paragraphs that followed the Assumptions
#    section.
# ------------------------------------------------------------------
$total = $d.Paragraphs.Count
$startPara = $null
for ($i = 1; $i -le $total; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13,[char]7) -eq "How to use") {
        $startPara = $i
        break
    }
}
if ($startPara -ne $null) {
    $startRange = $d.Paragraphs.Item($startPara).Range.Start
    $endRange = $d.Paragraphs.Item($total).Range.End
    $d.Range($startRange, $endRange).Delete() | Out-Null
}
